$wb = $excel.ActiveWorkbook

# --- "Order Sheet": clear Order Date / Item Description / Quantity on the
#     first data row of the OrdersPlaced table, then remove the remaining
#     data rows (3-17) so the table shrinks down to its header + 1 row. ---
$ws1 = $wb.Worksheets.Item("Order Sheet")
$ws1.Range("A2").ClearContents()
$ws1.Range("C2:D2").ClearContents()
$ws1.Rows("3:17").Delete()

# --- "Supplier Order Dates": clear out the two recorded supplier orders. ---
$ws3 = $wb.Worksheets.Item("Supplier Order Dates")
$ws3.Range("A2:B3").ClearContents()

# --- restore the on-screen selections for each sheet ---
$ws1.Select()
$ws1.Range("C10").Select()

$ws2 = $wb.Worksheets.Item("Order Helper")
$ws2.Select()
$ws2.Range("L13").Select()

$ws3.Select()
$ws3.Range("B3").Select()

$ws1.Select()
